$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell values
#    The fill order below is deliberately chosen so new shared-strings are
#    created in the same order as the target workbook (RunMode, Y, N,
#    Product, MotoX, "Apple Iphone 6", TestA, TestB).
# ---------------------------------------------------------------------------

# Old B1 ("Password") is no longer part of the sheet - drop it first.
$ws.Range("B1").Clear()

# Second little "TestB" table (rows 8-11) - column by column.
$ws.Range("A8").Value  = "RunMode"
$ws.Range("A9").Value  = "Y"
$ws.Range("A10").Value = "N"
$ws.Range("A11").Value = "Y"
$ws.Range("B8").Value  = "Product"
$ws.Range("B9").Value  = "MotoX"
$ws.Range("B10").Value = "Apple Iphone 6"
$ws.Range("B11").Value = "Akshay "

# First "TestA" table (rows 2-5).
$ws.Range("A2").Value = "RunMode"
$ws.Range("B2").Value = "UserName"
$ws.Range("C2").Value = "Password"
$ws.Range("A3").Value = "Omkar"
$ws.Range("B3").Value = "Omkar"
$ws.Range("C3").Value = "some"
$ws.Range("A4").Value = "Tejas"
$ws.Range("B4").Value = "Tejas"
$ws.Range("C4").Value = "more"
$ws.Range("A5").Value = "Akshay "
$ws.Range("B5").Value = "Akshay "
$ws.Range("C5").Value = "ssmsm"

# Table title labels, added last.
$ws.Range("A1").Value = "TestA"
$ws.Range("A7").Value = "TestB"

# ---------------------------------------------------------------------------
# 2. Formatting
# ---------------------------------------------------------------------------

# Thin black border around every populated cell of both tables.
$borderRange = $ws.Range("A1,A2:C5,A7:B7,A8:B11")
$borderRange.Borders.LineStyle = 1

# Yellow header rows for both tables.
$headerRange = $ws.Range("A2:C2,A8:B8")
$headerRange.Interior.Color = 65535

# Red "table title" cells (A1 / A7) use a plain (non-themed) Calibri font.
$titleRange = $ws.Range("A1,A7,B7")
$titleRange.Font.Name = "Calibri"
$titleRange.Interior.Color = 255

# Wrap text for the long product-name cell and size its row accordingly.
$ws.Range("B10").WrapText = $true
$ws.Rows(10).RowHeight = 28.8

# ---------------------------------------------------------------------------
# 3. Sheet/view level tweaks
# ---------------------------------------------------------------------------
$ws.Range("D5").Select()
$ws.PageSetup.Orientation = 1
